$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode subscript-three character used in row 22 price (e.g. 0.0₃0968)
$sub3 = [string][char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.774.04'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.104.33'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.99%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '389.57'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.57'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.544'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.43'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.73%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0861'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.595.52'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.69'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.86'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.099.44'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.989'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.89'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.835.50'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.50'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0' + $sub3 + '0968'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.05'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.92'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.14'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.20'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.15%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.28%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.16'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.110'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.37'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.60'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.07%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.91%  '

$ws.Range("B35").Value = 'OKB'

$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.26'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.16%  '

$ws.Range("B36").Value = 'VeChain'

$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0449'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.30%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.290'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.61%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.79%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.62'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.87'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.53%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.89%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.70'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.18'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.40%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.17%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.049.66'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.411.87'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.88%  '

$ws.Range("B51").Value = 'Mantle'

$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.902'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +15.71%  '
